$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting existing data right
$ws.Range("A1").EntireColumn.Insert()

# Populate the new column A with the image path for each menu row
$ws.Range("A1").Value = "images/menu/tuscan-grilled.jpg"
$ws.Range("A2").Value = "images/menu/tuscan-grilled.jpg"

# Set the width of the newly inserted column A (closest value the
# engine's pixel-grid ColumnWidth rounding can reproduce for 38.28515625)
$ws.Range("A1").EntireColumn.ColumnWidth = 37.5

# Update the active selection
$ws.Range("A2").Select()
